$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number that Excel would
# otherwise auto-convert from text to a numeric type; force them to
# stay text (matching the workbook's inlineStr/text convention) by
# briefly applying a text number format, then clearing the formatting
# so no residual style is left behind on the cell.
function Set-TextValue($ws, $ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).ClearFormats()
}

$ws.Range("D2").Value = "43.182.78"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "2.328.85"
$ws.Range("E3").Value = "  +1.25%  "
Set-TextValue $ws "D4" "1.00"
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue $ws "D5" "303.60"
$ws.Range("E5").Value = "  +1.20%  "
Set-TextValue $ws "D6" "97.87"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  -0.94%  "
Set-TextValue $ws "D8" "1.00"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -0.29%  "
Set-TextValue $ws "D10" "35.66"
$ws.Range("E10").Value = "  +0.03%  "
Set-TextValue $ws "D11" "19.29"
$ws.Range("E11").Value = "  +7.59%  "
Set-TextValue $ws "D12" "0.0791"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "2.684.41"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "2.331.20"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "43.095.19"
$ws.Range("E18").Value = "  +0.61%  "
Set-TextValue $ws "D19" "12.57"
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("E20").Value = "  -0.34%  "
Set-TextValue $ws "D21" "6.09"
$ws.Range("E21").Value = "  +0.82%  "
Set-TextValue $ws "D22" "67.97"
$ws.Range("E22").Value = "  +0.17%  "
Set-TextValue $ws "D23" "237.79"
$ws.Range("E23").Value = "  -0.90%  "
Set-TextValue $ws "D24" "2.20"
$ws.Range("E24").Value = "  +3.25%  "
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -1.34%  "
Set-TextValue $ws "D28" "2.37"
$ws.Range("E28").Value = "  +17.23%  "
Set-TextValue $ws "D29" "165.81"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("E30").Value = "  +0.70%  "
Set-TextValue $ws "D31" "33.17"
$ws.Range("E31").Value = "  +0.40%  "
Set-TextValue $ws "D32" "1.00"
$ws.Range("E32").Value = "  +0.05%  "
Set-TextValue $ws "D33" "18.10"
$ws.Range("E33").Value = "  +7.02%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E35").Value = "  -8.46%  "
$ws.Range("E36").Value = "  -1.88%  "
Set-TextValue $ws "D37" "0.0693"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("E39").Value = "  +1.85%  "
Set-TextValue $ws "D40" "1.76"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").Value = "2.001.85"
$ws.Range("E42").Value = "  -0.51%  "
Set-TextValue $ws "D43" "10.61"
$ws.Range("E43").Value = "  +4.42%  "
Set-TextValue $ws "D44" "0.0281"
$ws.Range("E44").Value = "  +0.36%  "
Set-TextValue $ws "D45" "18.34"
$ws.Range("E45").Value = "  +5.50%  "
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.554.23"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws "D49" "53.72"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").Value = "  -5.85%  "
Set-TextValue $ws "D51" "72.11"
$ws.Range("E51").Value = "  +0.11%  "
